# The deck currently applies the "Integral" theme (ppt/theme/theme2.xml) to
# the slide master / whole presentation, while ppt/theme/theme1.xml (unused
# by any slide, only referenced by the notes master) still holds the
# original default "Office Theme" palette. The commit swaps the two themes'
# contents so the presentation-visible theme becomes the plain default
# "Office Theme" color palette.
#
# The only theme surface exposed by the PowerPoint COM object model is the
# *currently active* ThemeColorScheme (Slide/SlideRange/CustomLayout), which
# always resolves to the theme part actually wired to the slide master
# (ppt/theme/theme2.xml in this deck). We drive the 12 theme colour slots to
# the default Office theme's RGB values via that API.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$tcs = $s.ThemeColorScheme

# ThemeColorSchemeIndex slot -> target RGB (packed 0xBBGGRR, as COM RGB()
# values are stored/read), taken from the stock "Office Theme" palette:
#  1 = dk1      000000
#  2 = lt1      FFFFFF
#  3 = dk2      44546A
#  4 = lt2      E7E6E6
#  5 = accent1  5B9BD5
#  6 = accent2  ED7D31
#  7 = accent3  A5A5A5
#  8 = accent4  FFC000
#  9 = accent5  4472C4
# 10 = accent6  70AD47
# 11 = hlink    0563C1
# 12 = folHlink 954F72
$tcs.Colors(1).RGB  = 0
$tcs.Colors(2).RGB  = 16777215
$tcs.Colors(3).RGB  = 6968388
$tcs.Colors(4).RGB  = 15132391
$tcs.Colors(5).RGB  = 13998939
$tcs.Colors(6).RGB  = 3243501
$tcs.Colors(7).RGB  = 10855845
$tcs.Colors(8).RGB  = 49407
$tcs.Colors(9).RGB  = 12874308
$tcs.Colors(10).RGB = 4697456
$tcs.Colors(11).RGB = 12673797
$tcs.Colors(12).RGB = 7491477
